$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, pushing the existing rows 44-109 down to 45-110.
$ws.Rows.Item(44).Insert()

# Populate the newly-inserted row 44 with the new weekly data point.
$ws.Cells.Item(44, 1).Value = 11
$ws.Cells.Item(44, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(44, 3).Value = 'Bíobío'
$ws.Cells.Item(44, 4).Value = 44495
$ws.Cells.Item(44, 5).Value = 8
$ws.Cells.Item(44, 6).Value = 'Fruta'
$ws.Cells.Item(44, 7).Value = 100108
$ws.Cells.Item(44, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(44, 9).Value = 100108005
$ws.Cells.Item(44, 10).Value = 'Piña'
$ws.Cells.Item(44, 11).Value = 'Caramelo'
$ws.Cells.Item(44, 12).Value = 'Segunda'
$ws.Cells.Item(44, 13).Value = 200
$ws.Cells.Item(44, 14).Value = 19000
$ws.Cells.Item(44, 15).Value = 20000
$ws.Cells.Item(44, 16).Value = 19500
$ws.Cells.Item(44, 17).Value = '$/caja 14 unidades'
$ws.Cells.Item(44, 18).Value = 'Ecuador'
$ws.Cells.Item(44, 19).Value = 1393
$ws.Cells.Item(44, 20).Value = 14
